$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 29   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/3/2022  Through  10/9/2022"

# --- Weekly crime-complaint figures (rows 14-30, columns C:N) ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = 14.285714285714
$ws.Range("I14").Value = 51
$ws.Range("J14").Value = 41
$ws.Range("K14").Value = 24.390243902439
$ws.Range("L14").Value = -3.77358490566
$ws.Range("M14").Value = -25
$ws.Range("N14").Value = -74.111675126903
$ws.Range("C15").Value = 4
$ws.Range("E15").Value = -20
$ws.Range("G15").Value = 19
$ws.Range("H15").Value = -5.263157894736
$ws.Range("I15").Value = 164
$ws.Range("J15").Value = 169
$ws.Range("K15").Value = -2.958579881656
$ws.Range("L15").Value = -13.684210526315
$ws.Range("M15").Value = 17.985611510791
$ws.Range("N15").Value = -62.124711316397
$ws.Range("C16").Value = 40
$ws.Range("E16").Value = 21.212121212121
$ws.Range("F16").Value = 158
$ws.Range("G16").Value = 118
$ws.Range("H16").Value = 33.898305084745
$ws.Range("I16").Value = 1511
$ws.Range("J16").Value = 1037
$ws.Range("K16").Value = 45.708775313404
$ws.Range("L16").Value = 20.590582601755
$ws.Range("M16").Value = -32.211754149843
$ws.Range("N16").Value = -85.867938645716
$ws.Range("C17").Value = 59
$ws.Range("D17").Value = 66
$ws.Range("E17").Value = -10.60606060606
$ws.Range("F17").Value = 231
$ws.Range("G17").Value = 279
$ws.Range("H17").Value = -17.204301075268
$ws.Range("I17").Value = 2671
$ws.Range("J17").Value = 2380
$ws.Range("K17").Value = 12.226890756302
$ws.Range("L17").Value = 23.030861354214
$ws.Range("M17").Value = 39.550679205851
$ws.Range("N17").Value = -49.953157204422
$ws.Range("C18").Value = 46
$ws.Range("D18").Value = 47
$ws.Range("E18").Value = -2.127659574468
$ws.Range("F18").Value = 180
$ws.Range("G18").Value = 198
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 1629
$ws.Range("J18").Value = 1342
$ws.Range("K18").Value = 21.385991058122
$ws.Range("L18").Value = 7.880794701986
$ws.Range("M18").Value = -39.442379182156
$ws.Range("N18").Value = -88.554767090564
$ws.Range("C19").Value = 121
$ws.Range("E19").Value = 0.833333333333
$ws.Range("F19").Value = 597
$ws.Range("G19").Value = 471
$ws.Range("H19").Value = 26.751592356687
$ws.Range("I19").Value = 5582
$ws.Range("J19").Value = 3813
$ws.Range("K19").Value = 46.393915552058
$ws.Range("L19").Value = 53.057307375925
$ws.Range("M19").Value = 32.431791221826
$ws.Range("N19").Value = -18.475244632685
$ws.Range("C20").Value = 42
$ws.Range("D20").Value = 30
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 151
$ws.Range("G20").Value = 121
$ws.Range("H20").Value = 24.793388429752
$ws.Range("I20").Value = 1428
$ws.Range("J20").Value = 986
$ws.Range("K20").Value = 44.827586206896
$ws.Range("L20").Value = 27.272727272727
$ws.Range("M20").Value = -6.176084099868
$ws.Range("N20").Value = -92.086011970738
$ws.Range("C21").Value = 313
$ws.Range("D21").Value = 303
$ws.Range("E21").Value = 3.300330033003
$ws.Range("F21").Value = 1343
$ws.Range("G21").Value = 1213
$ws.Range("H21").Value = 10.717230008244
$ws.Range("I21").Value = 13036
$ws.Range("J21").Value = 9768
$ws.Range("K21").Value = 33.456183456183
$ws.Range("L21").Value = 31.067765936054
$ws.Range("M21").Value = 2.027079909211
$ws.Range("N21").Value = -76.630873205098
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 15
$ws.Range("G22").Value = 13
$ws.Range("H22").Value = 15.384615384615
$ws.Range("I22").Value = 146
$ws.Range("J22").Value = 112
$ws.Range("K22").Value = 30.357142857142
$ws.Range("L22").Value = -3.311258278145
$ws.Range("M22").Value = -35.398230088495
$ws.Range("C23").Value = 16
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = 77.777777777777
$ws.Range("F23").Value = 41
$ws.Range("G23").Value = 41
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 431
$ws.Range("J23").Value = 376
$ws.Range("K23").Value = 14.627659574468
$ws.Range("L23").Value = 41.776315789473
$ws.Range("M23").Value = 61.423220973782
$ws.Range("C24").Value = 310
$ws.Range("D24").Value = 291
$ws.Range("E24").Value = 6.529209621993
$ws.Range("F24").Value = 1401
$ws.Range("G24").Value = 1067
$ws.Range("H24").Value = 31.302717900656
$ws.Range("I24").Value = 12578
$ws.Range("J24").Value = 8963
$ws.Range("K24").Value = 40.332477964967
$ws.Range("L24").Value = 35.553400150878
$ws.Range("M24").Value = 30.65337072816
$ws.Range("C25").Value = 105
$ws.Range("D25").Value = 114
$ws.Range("E25").Value = -7.894736842105
$ws.Range("F25").Value = 396
$ws.Range("G25").Value = 454
$ws.Range("H25").Value = -12.775330396475
$ws.Range("I25").Value = 4338
$ws.Range("J25").Value = 3760
$ws.Range("K25").Value = 15.372340425531
$ws.Range("L25").Value = 24.440619621342
$ws.Range("M25").Value = -17.669386980451
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 60
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 6.896551724137
$ws.Range("I26").Value = 260
$ws.Range("J26").Value = 265
$ws.Range("K26").Value = -1.88679245283
$ws.Range("L26").Value = -8.127208480565
$ws.Range("C27").Value = 13
$ws.Range("D27").Value = 21
$ws.Range("E27").Value = -38.095238095238
$ws.Range("F27").Value = 55
$ws.Range("G27").Value = 60
$ws.Range("H27").Value = -8.333333333333
$ws.Range("I27").Value = 540
$ws.Range("J27").Value = 482
$ws.Range("K27").Value = 12.033195020746
$ws.Range("L27").Value = 32.678132678132
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 17
$ws.Range("H28").Value = -11.764705882352
$ws.Range("I28").Value = 176
$ws.Range("J28").Value = 157
$ws.Range("K28").Value = 12.101910828025
$ws.Range("L28").Value = -24.137931034482
$ws.Range("M28").Value = -18.13953488372
$ws.Range("N28").Value = -72.107765451664
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 200
$ws.Range("F29").Value = 12
$ws.Range("G29").Value = 15
$ws.Range("H29").Value = -20
$ws.Range("I29").Value = 135
$ws.Range("J29").Value = 139
$ws.Range("K29").Value = -2.877697841726
$ws.Range("L29").Value = -23.728813559322
$ws.Range("M29").Value = -24.157303370786
$ws.Range("N29").Value = -75.364963503649
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = -50
$ws.Range("J30").Value = 52
$ws.Range("K30").Value = 78.846153846153
$ws.Range("L30").Value = 151.351351351351
